$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value, taken from the crypto-price refresh diff.
$updates = @(
    @{ Ref = 'D2'; Value = '27.155.92' },
    @{ Ref = 'E2'; Value = '  -1.97%  ' },
    @{ Ref = 'D3'; Value = '1.562.80' },
    @{ Ref = 'E4'; Value = '  +0.06%  ' },
    @{ Ref = 'D5'; Value = '206.74' },
    @{ Ref = 'E5'; Value = '  -0.27%  ' },
    @{ Ref = 'D6'; Value = '0.493' },
    @{ Ref = 'E6'; Value = '  -1.92%  ' },
    @{ Ref = 'E7'; Value = '  +0.08%  ' },
    @{ Ref = 'D8'; Value = '22.02' },
    @{ Ref = 'E8'; Value = '  -0.94%  ' },
    @{ Ref = 'E10'; Value = '  -0.17%  ' },
    @{ Ref = 'E11'; Value = '  -0.62%  ' },
    @{ Ref = 'D12'; Value = '1.785.29' },
    @{ Ref = 'E12'; Value = '  -1.71%  ' },
    @{ Ref = 'D13'; Value = '1.576.02' },
    @{ Ref = 'E13'; Value = '  -0.97%  ' },
    @{ Ref = 'E14'; Value = '  -2.67%  ' },
    @{ Ref = 'D15'; Value = '0.515' },
    @{ Ref = 'E15'; Value = '  -2.84%  ' },
    @{ Ref = 'D16'; Value = '63.02' },
    @{ Ref = 'E16'; Value = '  -0.77%  ' },
    @{ Ref = 'D17'; Value = '27.179.31' },
    @{ Ref = 'E17'; Value = '  -1.86%  ' },
    @{ Ref = 'E18'; Value = '  -1.18%  ' },
    @{ Ref = 'D19'; Value = '211.67' },
    @{ Ref = 'E19'; Value = '  -3.89%  ' },
    @{ Ref = 'D20'; Value = '7.21' },
    @{ Ref = 'E20'; Value = '  -1.89%  ' },
    @{ Ref = 'E21'; Value = '  +0.08%  ' },
    @{ Ref = 'D22'; Value = '4.11' },
    @{ Ref = 'E22'; Value = '  -0.89%  ' },
    @{ Ref = 'D23'; Value = '9.40' },
    @{ Ref = 'E23'; Value = '  -1.96%  ' },
    @{ Ref = 'E24'; Value = '  +0.26%  ' },
    @{ Ref = 'D25'; Value = '152.31' },
    @{ Ref = 'E25'; Value = '  -0.59%  ' },
    @{ Ref = 'D26'; Value = '6.62' },
    @{ Ref = 'E26'; Value = '  -3.76%  ' },
    @{ Ref = 'E27'; Value = '  -2.25%  ' },
    @{ Ref = 'E28'; Value = '  +0.08%  ' },
    @{ Ref = 'E29'; Value = '  -2.00%  ' },
    @{ Ref = 'D30'; Value = '1.14' },
    @{ Ref = 'E30'; Value = '  -0.94%  ' },
    @{ Ref = 'E31'; Value = '  -1.24%  ' },
    @{ Ref = 'D32'; Value = '3.17' },
    @{ Ref = 'E32'; Value = '  -1.80%  ' },
    @{ Ref = 'D33'; Value = '1.373.68' },
    @{ Ref = 'E33'; Value = '  +0.05%  ' },
    @{ Ref = 'E34'; Value = '  +0.43%  ' },
    @{ Ref = 'E36'; Value = '  -0.28%  ' },
    @{ Ref = 'D37'; Value = '0.943' },
    @{ Ref = 'E37'; Value = '  -3.68%  ' },
    @{ Ref = 'E38'; Value = '  -1.71%  ' },
    @{ Ref = 'D39'; Value = '0.521' },
    @{ Ref = 'E39'; Value = '  -3.24%  ' },
    @{ Ref = 'D40'; Value = '0.814' },
    @{ Ref = 'E40'; Value = '  -1.33%  ' },
    @{ Ref = 'E41'; Value = '  +0.14%  ' },
    @{ Ref = 'D42'; Value = '0.988' },
    @{ Ref = 'E42'; Value = '  +1.81%  ' },
    @{ Ref = 'D43'; Value = '1.80' },
    @{ Ref = 'E43'; Value = '  +3.46%  ' },
    @{ Ref = 'B44'; Value = 'MXToken' },
    @{ Ref = 'C44'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Ref = 'D44'; Value = '2.17' },
    @{ Ref = 'E44'; Value = '  +0.13%  ' },
    @{ Ref = 'B45'; Value = 'Aave' },
    @{ Ref = 'C45'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Ref = 'D45'; Value = '63.38' },
    @{ Ref = 'E45'; Value = '  -1.63%  ' },
    @{ Ref = 'E46'; Value = '  -0.94%  ' },
    @{ Ref = 'D47'; Value = '1.697.84' },
    @{ Ref = 'E47'; Value = '  -1.64%  ' },
    @{ Ref = 'D48'; Value = '85.43' },
    @{ Ref = 'E48'; Value = '  -2.85%  ' },
    @{ Ref = 'D49'; Value = '0.0₇0996' },
    @{ Ref = 'E49'; Value = '  -0.88%  ' },
    @{ Ref = 'E50'; Value = '  -1.03%  ' },
    @{ Ref = 'E51'; Value = '  +0.24%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    # Plain-decimal-looking strings (e.g. "206.74") get silently parsed into
    # numbers by Excel on assignment, which would drop the original text
    # representation. Force text storage for those, then drop the temporary
    # number-format override so the cell keeps its original (default) style.
    if ($u.Value -match '^[0-9]+\.[0-9]+$') {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
